$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column E to fit the new, longer "Metrics" summary ---
# (target stored width is 59.42578125 "characters"; this runtime quantizes
# ColumnWidth writes to 1/6-character steps, so 58.6666... is the input that
# lands on the closest reachable stored width, 59.5)
$ws.Columns.Item(5).ColumnWidth = 58.666666666666664

# --- Add new reference row (row 6), cloning the formatting of row 2 ---
$ws.Range("A2:E2").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)

$ws.Range("A6").Value = "Metrics for Measuring the Quality of Object-Oriented Software`n"
$ws.Range("B6").Value = "Gagandeep Singh"
$ws.Range("C6").Value = "9/1/2013"
$ws.Range("D6").Value = "http://delivery.acm.org.libproxy.auc.ca/10.1145/2510000/2507311/p66b-singh.pdf?ip=199.212.55.169&id=2507311&acc=ACTIVE%20SERVICE&key=FD0067F557510FFB%2E2E114FAB5F912086%2E4D4702B0C3E38B35%2E4D4702B0C3E38B35&CFID=939957675&CFTOKEN=22411919&__acm__=1495551123_66186ebadebfcf7a06c0c6f881edc6c7"
$ws.Range("E6").Value = "*LOC (lines of code) - understandability `n*CC (cyclomatic complexity) - complexity`n*CBO (coupling between objects) - efficiency, reuse, complexity`n*LCOM (lack of cohesion) - reuse, complexity`n*WMC (weighted methods per class) - maintainability, reuse`n*RFC (response for a class) - understandability, complexity`n*MI (maintainabiliyt index) - maintainability`n*NOC (number of children) - reuse, efficiency`n*DIT (depth of inheritance tree) - reuse, complexity "

$ws.Rows.Item(6).RowHeight = 135

# --- Update the selection to match the new active cell ---
$ws.Range("E6").Select() | Out-Null
